$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row 7 with the same pattern/style as the preceding rows.
$ws.Cells.Item(7, 1).Value = 42604.891469907408
$ws.Cells.Item(7, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(7, 2).Value = "Named"
$ws.Cells.Item(7, 3).Value = 2011
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = 0
